# Insert a new "PDF File Name" column (B) between the existing "Test Name"
# (A) and "Test Launch" (B->C) columns, then populate it with the PDF file
# names for each of the five A/B test rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts the old B:H block to C:I
# so the new (blank) column becomes B.
$ws.Columns("B:B").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# Header
$ws.Range("B1").Value = "PDF File Name"

# Data rows - PDF file names for each test (entered bottom-up so the
# shared-string table is populated in the same order as the source file)
$ws.Range("B6").Value = "5 - Homepage-Special-Offers-Carousel-Merchandising-Test-Results.pdf"
$ws.Range("B5").Value = "4 - Timber-Mountain-CTA-Copy-Test-Results.pdf"
$ws.Range("B4").Value = "3 - Timber-Mountain-Unified-Bundle-Flow-Checkout-Test-Results.pdf"
$ws.Range("B3").Value = "2 - Wild-Willy-AI-Planner-Trust-and-Adoption-AB-Test-Results.pdf"
$ws.Range("B2").Value = "1 - Locale-Aware-Experience-How-We-Boosted-International-Conversions-at-Timber-Mountain.pdf"

# Match formatting of the adjacent "Test Name" column (A) for the new
# PDF File Name column: header style + wrapped bold body style.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A2:A6").Copy()
$ws.Range("B2:B6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Update the active selection to reflect the new working location.
$ws.Range("A12").Select()
